$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = -7.833
$ws.Range("D3").Value = -7.757
$ws.Range("D5").Value = -7.961
$ws.Range("E7").Value = 13.045
$ws.Range("C9").Value = -11.899
$ws.Range("E9").Value = 12.659
$ws.Range("D11").Value = -7.699
$ws.Range("D12").Value = -7.644999999999999
$ws.Range("C13").Value = -12.201
$ws.Range("C16").Value = -12.439
$ws.Range("C18").Value = -12.362
$ws.Range("C20").Value = -12.32
$ws.Range("D21").Value = -7.897
$ws.Range("E21").Value = 13.084
